# Update gh-pages output generated at 456a3b4
# Applies numeric "want to go" (F column) bumps across sheets, plus a new
# "angela LIVE 2024" event row inserted into 演出 (appended) and 全部类型
# (inserted before the 2024.05.04 event, keeping chronological order).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - F column ("想去人数") updates
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 621
$ws1.Cells.Item(3, 6).Value = 610
$ws1.Cells.Item(5, 6).Value = 650
$ws1.Cells.Item(6, 6).Value = 798
$ws1.Cells.Item(7, 6).Value = 368
$ws1.Cells.Item(8, 6).Value = 571
$ws1.Cells.Item(9, 6).Value = 113
$ws1.Cells.Item(10, 6).Value = 1148
$ws1.Cells.Item(12, 6).Value = 349
$ws1.Cells.Item(13, 6).Value = 465
$ws1.Cells.Item(14, 6).Value = 150
$ws1.Cells.Item(17, 6).Value = 70
$ws1.Cells.Item(18, 6).Value = 530
$ws1.Cells.Item(19, 6).Value = 31
$ws1.Cells.Item(20, 6).Value = 534
$ws1.Cells.Item(21, 6).Value = 16
$ws1.Cells.Item(22, 6).Value = 518

# ---------------------------------------------------------------------
# Sheet "演出" (performances) - F column updates + new row 13
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 73
$ws2.Cells.Item(9, 6).Value = 199
$ws2.Cells.Item(11, 6).Value = 13

$r = 13
$ws2.Range("A12").Copy()
$ws2.Range("A13").PasteSpecial(-4122)
$ws2.Cells.Item($r, 1).Value = 12
$ws2.Cells.Item($r, 2).Value = "'2024.04.28"
$ws2.Cells.Item($r, 3).Value = "广州·「angela LIVE 2024」in  GUANGZHOU"
$ws2.Cells.Item($r, 4).Value = "奥体南路12号优托邦购物中心 疆进酒Omni Space GZ"
$ws2.Cells.Item($r, 5).Value = "2024.04.28 19:00-04.28 20:30"
$ws2.Cells.Item($r, 6).Value = 0
$ws2.Cells.Item($r, 7).Value = 480
$ws2.Cells.Item($r, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82041"
$ws2.Cells.Item($r, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) - F column updates + inserted row 31
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 621
$ws4.Cells.Item(7, 6).Value = 610
$ws4.Cells.Item(9, 6).Value = 650
$ws4.Cells.Item(10, 6).Value = 798
$ws4.Cells.Item(11, 6).Value = 368
$ws4.Cells.Item(12, 6).Value = 571
$ws4.Cells.Item(13, 6).Value = 113
$ws4.Cells.Item(14, 6).Value = 1148
$ws4.Cells.Item(18, 6).Value = 349
$ws4.Cells.Item(19, 6).Value = 465
$ws4.Cells.Item(21, 6).Value = 150
$ws4.Cells.Item(25, 6).Value = 70
$ws4.Cells.Item(26, 6).Value = 199
$ws4.Cells.Item(28, 6).Value = 530
$ws4.Cells.Item(29, 6).Value = 13

# Insert the new "angela LIVE 2024" event before the existing row 31
# (2024.05.04 黑塔利亚Only), shifting the remaining rows down by one.
$ws4.Rows.Item(31).Insert()

$r = 31
$ws4.Range("A30").Copy()
$ws4.Range("A31").PasteSpecial(-4122)
$ws4.Cells.Item($r, 1).Value = 30
$ws4.Cells.Item($r, 2).Value = "'2024.04.28"
$ws4.Cells.Item($r, 3).Value = "广州·「angela LIVE 2024」in  GUANGZHOU"
$ws4.Cells.Item($r, 4).Value = "奥体南路12号优托邦购物中心 疆进酒Omni Space GZ"
$ws4.Cells.Item($r, 5).Value = "2024.04.28 19:00-04.28 20:30"
$ws4.Cells.Item($r, 6).Value = 0
$ws4.Cells.Item($r, 7).Value = 480
$ws4.Cells.Item($r, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82041"
$ws4.Cells.Item($r, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg"

# Renumber the index column ("A") for the rows pushed down by the insert
# so it stays a plain 0..N sequential counter.
$ws4.Cells.Item(32, 1).Value = 31
$ws4.Cells.Item(33, 1).Value = 32
$ws4.Cells.Item(34, 1).Value = 33
$ws4.Cells.Item(35, 1).Value = 34
$ws4.Cells.Item(36, 1).Value = 35

# F-column values for the shifted rows (now 32-35) pick up the same
# "want to go" bumps applied elsewhere for these duplicated events.
$ws4.Cells.Item(32, 6).Value = 31
$ws4.Cells.Item(33, 6).Value = 534
$ws4.Cells.Item(34, 6).Value = 16
$ws4.Cells.Item(35, 6).Value = 518
